$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.847.04"
$ws.Range("E2").Value = "  -5.27%  "
$ws.Range("D3").Value = "3.209.56"
$ws.Range("E3").Value = "  -6.37%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'174.61"
$ws.Range("E5").Value = "  -7.21%  "
$ws.Range("D6").Value = "'513.60"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D7").Value = "'0.589"
$ws.Range("E7").Value = "  -4.93%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.211.35"
$ws.Range("E9").Value = "  -6.21%  "
$ws.Range("D10").Value = "'0.595"
$ws.Range("E10").Value = "  -6.89%  "
$ws.Range("D11").Value = "'52.33"
$ws.Range("E11").Value = "  -10.65%  "
$ws.Range("D12").Value = "'0.128"
$ws.Range("E12").Value = "  -6.84%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").Value = "'8.84"
$ws.Range("E14").Value = "  -6.95%  "
$ws.Range("D15").Value = "3.712.99"
$ws.Range("E15").Value = "  -6.17%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.115"
$ws.Range("E16").Value = "  -7.08%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.201.98"
$ws.Range("E17").Value = "  -6.03%  "
$ws.Range("D18").Value = "62.733.37"
$ws.Range("E18").Value = "  -4.88%  "
$ws.Range("D19").Value = "'17.07"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("D20").Value = "'10.88"
$ws.Range("E20").Value = "  -5.10%  "
$ws.Range("D21").Value = "'0.949"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").Value = "'362.53"
$ws.Range("E22").Value = "  -6.32%  "
$ws.Range("D23").Value = "'3.69"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "'11.02"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'79.72"
$ws.Range("E25").Value = "  -5.00%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'3.84"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.97"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("D28").Value = "'2.59"
$ws.Range("E28").Value = "  -5.49%  "
$ws.Range("D29").Value = "'11.16"
$ws.Range("E29").Value = "  -6.44%  "
$ws.Range("D30").Value = "'8.09"
$ws.Range("E30").Value = "  -6.93%  "
$ws.Range("D31").Value = "'648.37"
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("D32").Value = "'28.08"
$ws.Range("E32").Value = "  -7.06%  "
$ws.Range("D33").Value = "'6.22"
$ws.Range("E33").Value = "  -10.00%  "
$ws.Range("D34").Value = "'11.03"
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("D35").Value = "'0.103"
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("D36").Value = "'57.44"
$ws.Range("E36").Value = "  -7.88%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'36.24"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").Value = "'0.371"
$ws.Range("E39").Value = "  -4.87%  "
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").Value = "0.0₃0691"
$ws.Range("E41").Value = "  +7.13%  "
$ws.Range("D42").Value = "'0.121"
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("D43").Value = "2.841.42"
$ws.Range("E43").Value = "  -3.34%  "
$ws.Range("D44").Value = "'2.50"
$ws.Range("E44").Value = "  +2.46%  "
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "'0.0387"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.80"
$ws.Range("E47").Value = "  +6.59%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.56"
$ws.Range("E48").Value = "  -9.98%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'2.94"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'134.45"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.122"
$ws.Range("E51").Value = "  -4.43%  "
